# Updated cryptos list refresh (prices + 1h volume %) pulled in by the
# scheduled GitHub Actions job, plus a re-ranking swap of RenderToken and
# Mantle (rows 47/48 traded places).
#
# NOTE: several "Price" strings look numeric to Excel's smart-typing
# ("1.000", "29.407.13", ...) but must stay literal text, exactly like the
# rest of column D/E which are stored as plain strings (t="inlineStr").
# Forcing NumberFormat to "@" (Text) before the assignment stops Excel from
# reinterpreting them as numbers/dates; ClearFormats() afterwards drops the
# now-unneeded explicit number format again so the cell's style stays the
# untouched default, matching every other data cell in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Addr, $Val) {
    $rng = $ws.Range($Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.ClearFormats()
}

# row -> (new Price, new Volume(1h)); Price is $null where the diff left
# column D unchanged for that row.
$priceUpdates = @{
    2  = @("29.407.13", "  +0.21%  ")
    3  = @("1.876.03", "  +0.29%  ")
    4  = @("1.000", "  -0.06%  ")
    5  = @("0.7122", "  -0.43%  ")
    6  = @("241.84", "  +0.32%  ")
    7  = @($null, "  -0.02%  ")
    8  = @("0.3110", "  +0.71%  ")
    9  = @("0.07710", "  -2.42%  ")
    10 = @("25.43", "  +0.28%  ")
    11 = @("0.08381", "  +1.57%  ")
    12 = @("1.889.38", "  +0.74%  ")
    13 = @("5.251", "  +0.16%  ")
    14 = @("0.7163", "  -0.96%  ")
    15 = @("91.69", "  +1.01%  ")
    16 = @("29.426.79", "  +0.35%  ")
    17 = @("0.000008218", "  +4.98%  ")
    18 = @("5.981", "  +2.38%  ")
    19 = @("244.01", "  +0.09%  ")
    20 = @("2.136.83", "  +1.59%  ")
    21 = @("13.23", "  +0.08%  ")
    22 = @("0.9999", "  -0.06%  ")
    23 = @("7.923", "  -0.91%  ")
    24 = @("1.000", "  -0.06%  ")
    25 = @("0.1620", "  +1.28%  ")
    26 = @("163.60", "  +0.65%  ")
    27 = @("9.029", "  +0.51%  ")
    28 = @("18.60", "  +1.95%  ")
    29 = @($null, "  +0.79%  ")
    30 = @("4.420", "  +0.98%  ")
    31 = @("1.297", "  -3.80%  ")
    32 = @("4.317", "  +5.28%  ")
    33 = @("0.05226", "  +0.73%  ")
    34 = @("1.926", "  -0.74%  ")
    35 = @("0.7748", "  +6.98%  ")
    36 = @("1.175", "  -1.00%  ")
    37 = @("2.682", "  +0.31%  ")
    38 = @("0.01866", "  +0.53%  ")
    39 = @("2.722", "  +0.91%  ")
    40 = @("1.166.64", "  -0.48%  ")
    41 = @("6.419", "  +4.81%  ")
    42 = @("73.51", "  +1.27%  ")
    43 = @("0.8912", "  -1.48%  ")
    44 = @("104.54", "  +2.62%  ")
    45 = @("1.000", "  -0.08%  ")
    46 = @("2.033.39", "  +1.16%  ")
    49 = @("9.411", "  +1.53%  ")
    50 = @("0.4309", "  +0.71%  ")
    51 = @("7.074", "  +0.67%  ")
}

foreach ($row in $priceUpdates.Keys) {
    $pair = $priceUpdates[$row]
    $price = $pair[0]
    $volume = $pair[1]
    if ($null -ne $price) {
        Set-TextValue "D$row" $price
    }
    Set-TextValue "E$row" $volume
}

# Rows 47/48 also swapped rank (RenderToken <-> Mantle), each bringing its
# own new Coin/Link/Price/Volume values with it.
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D47" "0.5205"
Set-TextValue "E47" "  -1.48%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "1.796"
Set-TextValue "E48" "  +0.55%  "
